$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the company name and related id in row 2
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "societé B"

# Update selection to B2 (mission table selection for loading plan display)
$ws.Range("B2").Select()
